$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value2 = 3822.111
$ws.Range("I74").Value2 = 3299.8
$ws.Range("J74").Value2 = 4475
$ws.Range("K74").Value2 = 3299.8
$ws.Range("L74").Value2 = 4475
$ws.Range("M74").Value2 = -2363.8
$ws.Range("N74").Value2 = -6347
$ws.Range("H77").Value2 = 3822.111
$ws.Range("I77").Value2 = 3299.8
$ws.Range("J77").Value2 = 4475
$ws.Range("K77").Value2 = 16499
$ws.Range("L77").Value2 = 22375
$ws.Range("M77").Value2 = -11819
$ws.Range("N77").Value2 = -31735
$ws.Range("H86").Value2 = 3698.4546
$ws.Range("I86").Value2 = 2973.25
$ws.Range("J86").Value2 = 4112.857
$ws.Range("K86").Value2 = 2973.25
$ws.Range("L86").Value2 = 4112.857
$ws.Range("M86").Value2 = -1850.25
$ws.Range("N86").Value2 = -6358.857
$ws.Range("H89").Value2 = 3698.4546
$ws.Range("I89").Value2 = 2973.25
$ws.Range("J89").Value2 = 4112.857
$ws.Range("K89").Value2 = 14866.25
$ws.Range("L89").Value2 = 20564.285
$ws.Range("M89").Value2 = -9250.25
$ws.Range("N89").Value2 = -31796.285
$ws.Range("H99").Value2 = 12709.125
$ws.Range("I99").Value2 = 20135.6
$ws.Range("J99").Value2 = 331.66666
$ws.Range("K99").Value2 = 60406.8
$ws.Range("L99").Value2 = 994.9999799999999
$ws.Range("M99").Value2 = -58908.8
$ws.Range("N99").Value2 = -3990.99998
$ws.Range("H116").Value2 = 1671
$ws.Range("I116").Value2 = 1574.1
$ws.Range("J116").Value2 = 1832.5
$ws.Range("K116").Value2 = 1574.1
$ws.Range("L116").Value2 = 1832.5
$ws.Range("M116").Value2 = 1867.9
$ws.Range("N116").Value2 = -8716.5
$ws.Range("H132").Value2 = 4907617
$ws.Range("I132").Value2 = 5561390
$ws.Range("J132").Value2 = 4317.6665
$ws.Range("K132").Value2 = 16684170
$ws.Range("L132").Value2 = 12952.9995
$ws.Range("M132").Value2 = -16681640
$ws.Range("N132").Value2 = -18012.9995
$ws.Range("H137").Value2 = 1656.0869
$ws.Range("I137").Value2 = 1694.8334
$ws.Range("K137").Value2 = 5084.5002
$ws.Range("M137").Value2 = -2534.5002
$ws.Range("H138").Value2 = 4302.674
$ws.Range("I138").Value2 = 3682.7
$ws.Range("J138").Value2 = 4474.8887
$ws.Range("K138").Value2 = 11048.1
$ws.Range("L138").Value2 = 13424.6661
$ws.Range("M138").Value2 = -5908.099999999999
$ws.Range("N138").Value2 = -23704.6661

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 37875.38
$ws.Range("I32").Value2 = 16631.174
$ws.Range("J32").Value2 = 159546.73
$ws.Range("K32").Value2 = 16631.174
$ws.Range("L32").Value2 = 159546.73
$ws.Range("M32").Value2 = -16344.174
$ws.Range("N32").Value2 = -160120.73
$ws.Range("H45").Value2 = 59979.824
$ws.Range("I45").Value2 = 91637.55
$ws.Range("J45").Value2 = 1940.6666
$ws.Range("K45").Value2 = 91637.55
$ws.Range("L45").Value2 = 1940.6666
$ws.Range("M45").Value2 = -91260.55
$ws.Range("N45").Value2 = -2694.6666
$ws.Range("H55").Value2 = 12122.5
$ws.Range("J55").Value2 = 12122.5
$ws.Range("L55").Value2 = 12122.5
$ws.Range("N55").Value2 = -12752.5
$ws.Range("H74").Value2 = 1855
$ws.Range("I74").Value2 = 1920.6
$ws.Range("J74").Value2 = 1773
$ws.Range("K74").Value2 = 1920.6
$ws.Range("L74").Value2 = 1773
$ws.Range("M74").Value2 = -1046.6
$ws.Range("N74").Value2 = -3521
$ws.Range("H77").Value2 = 1855
$ws.Range("I77").Value2 = 1920.6
$ws.Range("J77").Value2 = 1773
$ws.Range("K77").Value2 = 9603
$ws.Range("L77").Value2 = 8865
$ws.Range("M77").Value2 = -5235
$ws.Range("N77").Value2 = -17601
$ws.Range("H110").Value2 = 45551230
$ws.Range("I110").Value2 = 55673440
$ws.Range("J110").Value2 = 1312.5
$ws.Range("K110").Value2 = 55673440
$ws.Range("L110").Value2 = 1312.5
$ws.Range("M110").Value2 = -55671395
$ws.Range("N110").Value2 = -5402.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value2 = 2023.3334
$ws.Range("I99").Value2 = 1917.1428
$ws.Range("K99").Value2 = 1917.1428
$ws.Range("M99").Value2 = -419.1428000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 45864.312
$ws.Range("I31").Value2 = 1916.2307
$ws.Range("J31").Value2 = 71833.63
$ws.Range("K31").Value2 = 1916.2307
$ws.Range("L31").Value2 = 71833.63
$ws.Range("M31").Value2 = -1621.2307
$ws.Range("N31").Value2 = -72423.63
$ws.Range("H34").Value2 = 45864.312
$ws.Range("I34").Value2 = 1916.2307
$ws.Range("J34").Value2 = 71833.63
$ws.Range("K34").Value2 = 1916.2307
$ws.Range("L34").Value2 = 71833.63
$ws.Range("M34").Value2 = -1714.2307
$ws.Range("N34").Value2 = -72237.63
$ws.Range("H134").Value2 = 1172.7646
$ws.Range("I134").Value2 = 727.6316
$ws.Range("J134").Value2 = 1736.6
$ws.Range("K134").Value2 = 2182.8948
$ws.Range("L134").Value2 = 5209.799999999999
$ws.Range("M134").Value2 = 352.1052
$ws.Range("N134").Value2 = -10279.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value2 = 969
$ws.Range("I68").Value2 = 615
$ws.Range("J68").Value2 = 1110.6
$ws.Range("K68").Value2 = 1845
$ws.Range("L68").Value2 = 3331.8
$ws.Range("M68").Value2 = -1034
$ws.Range("N68").Value2 = -4953.799999999999
$ws.Range("H71").Value2 = 969
$ws.Range("I71").Value2 = 615
$ws.Range("J71").Value2 = 1110.6
$ws.Range("K71").Value2 = 5535
$ws.Range("L71").Value2 = 9995.4
$ws.Range("M71").Value2 = -1479
$ws.Range("N71").Value2 = -18107.4
$ws.Range("H88").Value2 = 1750
$ws.Range("J88").Value2 = 1750
$ws.Range("L88").Value2 = 5250
$ws.Range("N88").Value2 = -6106
$ws.Range("H91").Value2 = 1750
$ws.Range("J91").Value2 = 1750
$ws.Range("L91").Value2 = 5250
$ws.Range("N91").Value2 = -8214
$ws.Range("H113").Value2 = 1145.2354
$ws.Range("I113").Value2 = 1958.4286
$ws.Range("J113").Value2 = 576
$ws.Range("K113").Value2 = 5875.2858
$ws.Range("L113").Value2 = 1728
$ws.Range("M113").Value2 = -3705.2858
$ws.Range("N113").Value2 = -6068
$ws.Range("H131").Value2 = 704281.1
$ws.Range("J131").Value2 = 780052.9399999999
$ws.Range("L131").Value2 = 2340158.82
$ws.Range("N131").Value2 = -2350238.82

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value2 = 0
$ws.Range("I3").Value2 = 0
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 0
$ws.Range("L3").Value2 = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H126").Value2 = 3696
$ws.Range("I126").Value2 = 3120
$ws.Range("J126").Value2 = 6000
$ws.Range("K126").Value2 = 9360
$ws.Range("L126").Value2 = 18000
$ws.Range("M126").Value2 = -6890
$ws.Range("N126").Value2 = -22940

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value2 = 1599.1471
$ws.Range("I68").Value2 = 1566.8182
$ws.Range("J68").Value2 = 1614.6086
$ws.Range("K68").Value2 = 1566.8182
$ws.Range("L68").Value2 = 1614.6086
$ws.Range("M68").Value2 = -817.8181999999999
$ws.Range("N68").Value2 = -3112.6086
$ws.Range("H71").Value2 = 1599.1471
$ws.Range("I71").Value2 = 1566.8182
$ws.Range("J71").Value2 = 1614.6086
$ws.Range("K71").Value2 = 7834.090999999999
$ws.Range("L71").Value2 = 8073.043
$ws.Range("M71").Value2 = -4090.090999999999
$ws.Range("N71").Value2 = -15561.043
